# Apply the 2022-06-07 data refresh to the "Fonds de solidarite" workbook.
# For each impacted row, update the "nombre_aides" (column C) and
# "montant_total" (column E) values to their new totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new C (nombre_aides) / E (montant_total) values
$updates = @(
    @{ Row = 20;  C = 1593;   E = 53358905   },
    @{ Row = 26;  C = 33938;  E = 203949247  },
    @{ Row = 44;  C = 10556;  E = 42605161   },
    @{ Row = 74;  C = 951;    E = 4274646    },
    @{ Row = 77;  C = 4544;   E = 8328779    },
    @{ Row = 92;  C = 409150; E = 1595480068 },
    @{ Row = 94;  C = 94199;  E = 917978553  },
    @{ Row = 95;  C = 50772;  E = 932793689  },
    @{ Row = 96;  C = 17294;  E = 794330219  },
    @{ Row = 97;  C = 2157;   E = 214111625  },
    @{ Row = 110; C = 397;    E = 16698678   },
    @{ Row = 141; C = 80474;  E = 280718664  },
    @{ Row = 174; C = 226092; E = 900594909  }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
